$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5
$ws.Cells.Item(5, 2).Value = 5265388
$ws.Cells.Item(5, 7).Value = "UD Oliveirense"
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 10).Value = "D"
$ws.Cells.Item(5, 11).Value = 2.2
$ws.Cells.Item(5, 12).Value = 3.2
$ws.Cells.Item(5, 13).Value = 2.875
$ws.Cells.Item(5, 14).Value = 1.909
$ws.Cells.Item(5, 16).Value = 3.6
$ws.Cells.Item(5, 18).Value = 1.9
$ws.Cells.Item(5, 19).Value = 1.9
$ws.Cells.Item(5, 21).Value = 1.975
$ws.Cells.Item(5, 22).Value = 1.825
$ws.Cells.Item(5, 24).Value = 2.5
$ws.Cells.Item(5, 25).Value = -1
$ws.Cells.Item(5, 27).Value = 0.8999999999999999
$ws.Cells.Item(5, 29).Value = 0.825

# Row 6
$ws.Cells.Item(6, 2).Value = 5266496
$ws.Cells.Item(6, 7).Value = "Belenenses"
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 10).Value = "A"
$ws.Cells.Item(6, 11).Value = 1.833
$ws.Cells.Item(6, 12).Value = 3.4
$ws.Cells.Item(6, 13).Value = 3.6
$ws.Cells.Item(6, 14).Value = 1.75
$ws.Cells.Item(6, 16).Value = 4.2
$ws.Cells.Item(6, 18).Value = 1.8
$ws.Cells.Item(6, 19).Value = 2
$ws.Cells.Item(6, 21).Value = 1.875
$ws.Cells.Item(6, 22).Value = 1.925
$ws.Cells.Item(6, 24).Value = -1
$ws.Cells.Item(6, 25).Value = 3.2
$ws.Cells.Item(6, 27).Value = 1
$ws.Cells.Item(6, 29).Value = 0.925

# Row 41
$ws.Cells.Item(41, 2).Value = 5266489
$ws.Cells.Item(41, 6).Value = "FC Porto B"
$ws.Cells.Item(41, 7).Value = "Belenenses"
$ws.Cells.Item(41, 8).Value = 1
$ws.Cells.Item(41, 9).Value = 2
$ws.Cells.Item(41, 10).Value = "A"
$ws.Cells.Item(41, 11).Value = 1.75
$ws.Cells.Item(41, 12).Value = 3.4
$ws.Cells.Item(41, 13).Value = 4.2
$ws.Cells.Item(41, 14).Value = 1.6
$ws.Cells.Item(41, 15).Value = 3.5
$ws.Cells.Item(41, 16).Value = 5
$ws.Cells.Item(41, 17).Value = -0.75
$ws.Cells.Item(41, 18).Value = 1.775
$ws.Cells.Item(41, 19).Value = 2.025
$ws.Cells.Item(41, 20).Value = 2.5
$ws.Cells.Item(41, 21).Value = 1.925
$ws.Cells.Item(41, 22).Value = 1.875
$ws.Cells.Item(41, 23).Value = -1
$ws.Cells.Item(41, 25).Value = 4
$ws.Cells.Item(41, 26).Value = -1
$ws.Cells.Item(41, 27).Value = 1.025
$ws.Cells.Item(41, 28).Value = 0.925

# Row 42
$ws.Cells.Item(42, 2).Value = 5265414
$ws.Cells.Item(42, 6).Value = "Nacional"
$ws.Cells.Item(42, 7).Value = "SCU Torreense"
$ws.Cells.Item(42, 8).Value = 3
$ws.Cells.Item(42, 9).Value = 1
$ws.Cells.Item(42, 10).Value = "H"
$ws.Cells.Item(42, 11).Value = 2.2
$ws.Cells.Item(42, 12).Value = 3.1
$ws.Cells.Item(42, 13).Value = 3
$ws.Cells.Item(42, 14).Value = 2.45
$ws.Cells.Item(42, 15).Value = 3
$ws.Cells.Item(42, 16).Value = 2.7
$ws.Cells.Item(42, 17).Value = 0
$ws.Cells.Item(42, 18).Value = 1.825
$ws.Cells.Item(42, 19).Value = 1.975
$ws.Cells.Item(42, 20).Value = 2.25
$ws.Cells.Item(42, 21).Value = 2.025
$ws.Cells.Item(42, 22).Value = 1.775
$ws.Cells.Item(42, 23).Value = 1.45
$ws.Cells.Item(42, 25).Value = -1
$ws.Cells.Item(42, 26).Value = 0.825
$ws.Cells.Item(42, 27).Value = -1
$ws.Cells.Item(42, 28).Value = 1.025

# Row 44
$ws.Cells.Item(44, 2).Value = 5263141
$ws.Cells.Item(44, 6).Value = "Trofense"
$ws.Cells.Item(44, 7).Value = "Leixoes"
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = "A"
$ws.Cells.Item(44, 11).Value = 2.6
$ws.Cells.Item(44, 12).Value = 3.1
$ws.Cells.Item(44, 14).Value = 2.55
$ws.Cells.Item(44, 15).Value = 3
$ws.Cells.Item(44, 16).Value = 2.55
$ws.Cells.Item(44, 18).Value = 1.9
$ws.Cells.Item(44, 19).Value = 1.9
$ws.Cells.Item(44, 20).Value = 2.25
$ws.Cells.Item(44, 21).Value = 1.9
$ws.Cells.Item(44, 22).Value = 1.9
$ws.Cells.Item(44, 23).Value = -1
$ws.Cells.Item(44, 25).Value = 1.55
$ws.Cells.Item(44, 26).Value = -1
$ws.Cells.Item(44, 27).Value = 0.8999999999999999
$ws.Cells.Item(44, 28).Value = -1
$ws.Cells.Item(44, 29).Value = 0.8999999999999999

# Row 45
$ws.Cells.Item(45, 2).Value = 5265413
$ws.Cells.Item(45, 6).Value = "Benfica B"
$ws.Cells.Item(45, 7).Value = "Tondela"
$ws.Cells.Item(45, 8).Value = 2
$ws.Cells.Item(45, 10).Value = "H"
$ws.Cells.Item(45, 11).Value = 2.55
$ws.Cells.Item(45, 12).Value = 3.2
$ws.Cells.Item(45, 14).Value = 2.375
$ws.Cells.Item(45, 15).Value = 3.2
$ws.Cells.Item(45, 16).Value = 2.75
$ws.Cells.Item(45, 18).Value = 1.75
$ws.Cells.Item(45, 19).Value = 2.05
$ws.Cells.Item(45, 20).Value = 2.5
$ws.Cells.Item(45, 21).Value = 1.825
$ws.Cells.Item(45, 22).Value = 1.975
$ws.Cells.Item(45, 23).Value = 1.375
$ws.Cells.Item(45, 25).Value = -1
$ws.Cells.Item(45, 26).Value = 0.75
$ws.Cells.Item(45, 27).Value = -1
$ws.Cells.Item(45, 28).Value = 0.825
$ws.Cells.Item(45, 29).Value = -1

# Row 63
$ws.Cells.Item(63, 2).Value = 5265427
$ws.Cells.Item(63, 6).Value = "Benfica B"
$ws.Cells.Item(63, 9).Value = 2
$ws.Cells.Item(63, 10).Value = "A"
$ws.Cells.Item(63, 11).Value = 2.15
$ws.Cells.Item(63, 12).Value = 3.2
$ws.Cells.Item(63, 13).Value = 3.2
$ws.Cells.Item(63, 14).Value = 2.5
$ws.Cells.Item(63, 15).Value = 3.3
$ws.Cells.Item(63, 16).Value = 2.4
$ws.Cells.Item(63, 17).Value = 0
$ws.Cells.Item(63, 18).Value = 1.95
$ws.Cells.Item(63, 19).Value = 1.85
$ws.Cells.Item(63, 21).Value = 1.9
$ws.Cells.Item(63, 22).Value = 1.9
$ws.Cells.Item(63, 23).Value = -1
$ws.Cells.Item(63, 25).Value = 1.4
$ws.Cells.Item(63, 26).Value = -1
$ws.Cells.Item(63, 27).Value = 0.8500000000000001
$ws.Cells.Item(63, 28).Value = 0.8999999999999999
$ws.Cells.Item(63, 29).Value = -1

# Row 64
$ws.Cells.Item(64, 2).Value = 5265428
$ws.Cells.Item(64, 6).Value = "SC Farense"
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = "H"
$ws.Cells.Item(64, 11).Value = 1.909
$ws.Cells.Item(64, 12).Value = 3.25
$ws.Cells.Item(64, 13).Value = 3.75
$ws.Cells.Item(64, 14).Value = 1.833
$ws.Cells.Item(64, 15).Value = 3.4
$ws.Cells.Item(64, 16).Value = 4
$ws.Cells.Item(64, 17).Value = -0.5
$ws.Cells.Item(64, 18).Value = 1.825
$ws.Cells.Item(64, 19).Value = 1.975
$ws.Cells.Item(64, 21).Value = 1.975
$ws.Cells.Item(64, 22).Value = 1.825
$ws.Cells.Item(64, 23).Value = 0.833
$ws.Cells.Item(64, 25).Value = -1
$ws.Cells.Item(64, 26).Value = 0.825
$ws.Cells.Item(64, 27).Value = -1
$ws.Cells.Item(64, 28).Value = -1
$ws.Cells.Item(64, 29).Value = 0.825

# Row 74
$ws.Cells.Item(74, 2).Value = 5265438
$ws.Cells.Item(74, 6).Value = "Penafiel"
$ws.Cells.Item(74, 7).Value = "UD Oliveirense"
$ws.Cells.Item(74, 8).Value = 1
$ws.Cells.Item(74, 9).Value = 1
$ws.Cells.Item(74, 10).Value = "D"
$ws.Cells.Item(74, 11).Value = 2.05
$ws.Cells.Item(74, 12).Value = 3.2
$ws.Cells.Item(74, 13).Value = 3.4
$ws.Cells.Item(74, 14).Value = 1.909
$ws.Cells.Item(74, 16).Value = 3.6
$ws.Cells.Item(74, 17).Value = -0.5
$ws.Cells.Item(74, 18).Value = 1.95
$ws.Cells.Item(74, 19).Value = 1.85
$ws.Cells.Item(74, 20).Value = 2.5
$ws.Cells.Item(74, 21).Value = 2
$ws.Cells.Item(74, 22).Value = 1.8
$ws.Cells.Item(74, 24).Value = 2.4
$ws.Cells.Item(74, 25).Value = -1
$ws.Cells.Item(74, 27).Value = 0.8500000000000001
$ws.Cells.Item(74, 29).Value = 0.8

# Row 75
$ws.Cells.Item(75, 2).Value = 5265440
$ws.Cells.Item(75, 6).Value = "Benfica B"
$ws.Cells.Item(75, 7).Value = "Vilafranquense"
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 2
$ws.Cells.Item(75, 10).Value = "A"
$ws.Cells.Item(75, 11).Value = 3
$ws.Cells.Item(75, 12).Value = 3.3
$ws.Cells.Item(75, 13).Value = 2.25
$ws.Cells.Item(75, 14).Value = 2.15
$ws.Cells.Item(75, 16).Value = 3
$ws.Cells.Item(75, 17).Value = -0.25
$ws.Cells.Item(75, 18).Value = 1.9
$ws.Cells.Item(75, 19).Value = 1.9
$ws.Cells.Item(75, 20).Value = 2.75
$ws.Cells.Item(75, 21).Value = 1.875
$ws.Cells.Item(75, 22).Value = 1.925
$ws.Cells.Item(75, 24).Value = -1
$ws.Cells.Item(75, 25).Value = 2
$ws.Cells.Item(75, 27).Value = 0.8999999999999999
$ws.Cells.Item(75, 29).Value = 0.925

# Row 81
$ws.Cells.Item(81, 2).Value = 5266482
$ws.Cells.Item(81, 6).Value = "Feirense"
$ws.Cells.Item(81, 7).Value = "Belenenses"
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = "A"
$ws.Cells.Item(81, 11).Value = 1.8
$ws.Cells.Item(81, 12).Value = 3.4
$ws.Cells.Item(81, 13).Value = 4
$ws.Cells.Item(81, 14).Value = 1.75
$ws.Cells.Item(81, 15).Value = 3.5
$ws.Cells.Item(81, 16).Value = 4.2
$ws.Cells.Item(81, 17).Value = -0.75
$ws.Cells.Item(81, 18).Value = 2
$ws.Cells.Item(81, 19).Value = 1.8
$ws.Cells.Item(81, 20).Value = 2.25
$ws.Cells.Item(81, 21).Value = 1.775
$ws.Cells.Item(81, 22).Value = 2.025
$ws.Cells.Item(81, 24).Value = -1
$ws.Cells.Item(81, 25).Value = 3.2
$ws.Cells.Item(81, 26).Value = -1
$ws.Cells.Item(81, 27).Value = 0.8
$ws.Cells.Item(81, 29).Value = 1.025

# Row 82
$ws.Cells.Item(82, 2).Value = 5265435
$ws.Cells.Item(82, 6).Value = "FC Porto B"
$ws.Cells.Item(82, 7).Value = "CD Mafra"
$ws.Cells.Item(82, 8).Value = 1
$ws.Cells.Item(82, 10).Value = "D"
$ws.Cells.Item(82, 11).Value = 2.05
$ws.Cells.Item(82, 12).Value = 3.3
$ws.Cells.Item(82, 13).Value = 3.3
$ws.Cells.Item(82, 14).Value = 2
$ws.Cells.Item(82, 15).Value = 3.4
$ws.Cells.Item(82, 16).Value = 3.4
$ws.Cells.Item(82, 17).Value = -0.25
$ws.Cells.Item(82, 18).Value = 1.775
$ws.Cells.Item(82, 19).Value = 2.025
$ws.Cells.Item(82, 20).Value = 2.5
$ws.Cells.Item(82, 21).Value = 1.9
$ws.Cells.Item(82, 22).Value = 1.9
$ws.Cells.Item(82, 24).Value = 2.4
$ws.Cells.Item(82, 25).Value = -1
$ws.Cells.Item(82, 26).Value = -0.5
$ws.Cells.Item(82, 27).Value = 0.5125
$ws.Cells.Item(82, 29).Value = 0.8999999999999999

# Row 90
$ws.Cells.Item(90, 2).Value = 5265446
$ws.Cells.Item(90, 6).Value = "Leixoes"
$ws.Cells.Item(90, 7).Value = "Estrela"
$ws.Cells.Item(90, 8).Value = 1
$ws.Cells.Item(90, 9).Value = 1
$ws.Cells.Item(90, 10).Value = "D"
$ws.Cells.Item(90, 11).Value = 2.875
$ws.Cells.Item(90, 12).Value = 3.25
$ws.Cells.Item(90, 13).Value = 2.25
$ws.Cells.Item(90, 14).Value = 3.25
$ws.Cells.Item(90, 15).Value = 3.2
$ws.Cells.Item(90, 16).Value = 2.05
$ws.Cells.Item(90, 18).Value = 2
$ws.Cells.Item(90, 19).Value = 1.85
$ws.Cells.Item(90, 21).Value = 1.85
$ws.Cells.Item(90, 22).Value = 2
$ws.Cells.Item(90, 24).Value = 2.2
$ws.Cells.Item(90, 25).Value = -1
$ws.Cells.Item(90, 26).Value = 0.5
$ws.Cells.Item(90, 27).Value = -0.5
$ws.Cells.Item(90, 29).Value = 0.5

# Row 91
$ws.Cells.Item(91, 2).Value = 5266479
$ws.Cells.Item(91, 6).Value = "Belenenses"
$ws.Cells.Item(91, 7).Value = "Tondela"
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 2
$ws.Cells.Item(91, 10).Value = "A"
$ws.Cells.Item(91, 11).Value = 2.5
$ws.Cells.Item(91, 12).Value = 3.1
$ws.Cells.Item(91, 13).Value = 2.625
$ws.Cells.Item(91, 14).Value = 3
$ws.Cells.Item(91, 15).Value = 3
$ws.Cells.Item(91, 16).Value = 2.3
$ws.Cells.Item(91, 18).Value = 1.775
$ws.Cells.Item(91, 19).Value = 2.025
$ws.Cells.Item(91, 21).Value = 2
$ws.Cells.Item(91, 22).Value = 1.8
$ws.Cells.Item(91, 24).Value = -1
$ws.Cells.Item(91, 25).Value = 1.3
$ws.Cells.Item(91, 26).Value = -1
$ws.Cells.Item(91, 27).Value = 1.025
$ws.Cells.Item(91, 29).Value = 0.4

# Row 116
$ws.Cells.Item(116, 2).Value = 5265464
$ws.Cells.Item(116, 6).Value = "SCU Torreense"
$ws.Cells.Item(116, 7).Value = "Sporting Covilha"
$ws.Cells.Item(116, 8).Value = 3
$ws.Cells.Item(116, 11).Value = 1.8
$ws.Cells.Item(116, 12).Value = 3.4
$ws.Cells.Item(116, 13).Value = 4
$ws.Cells.Item(116, 14).Value = 1.666
$ws.Cells.Item(116, 15).Value = 3.6
$ws.Cells.Item(116, 16).Value = 4.333
$ws.Cells.Item(116, 17).Value = -0.75
$ws.Cells.Item(116, 18).Value = 1.975
$ws.Cells.Item(116, 19).Value = 1.825
$ws.Cells.Item(116, 21).Value = 1.95
$ws.Cells.Item(116, 22).Value = 1.85
$ws.Cells.Item(116, 23).Value = 0.6659999999999999
$ws.Cells.Item(116, 26).Value = 0.9750000000000001
$ws.Cells.Item(116, 28).Value = 0.95

# Row 117
$ws.Cells.Item(117, 2).Value = 5263148
$ws.Cells.Item(117, 6).Value = "FC Porto B"
$ws.Cells.Item(117, 7).Value = "Trofense"
$ws.Cells.Item(117, 8).Value = 4
$ws.Cells.Item(117, 9).Value = 2
$ws.Cells.Item(117, 11).Value = 1.615
$ws.Cells.Item(117, 12).Value = 3.8
$ws.Cells.Item(117, 13).Value = 4.75
$ws.Cells.Item(117, 16).Value = 4.5
$ws.Cells.Item(117, 18).Value = 1.875
$ws.Cells.Item(117, 19).Value = 1.925
$ws.Cells.Item(117, 20).Value = 2.5
$ws.Cells.Item(117, 21).Value = 1.85
$ws.Cells.Item(117, 22).Value = 1.95
$ws.Cells.Item(117, 26).Value = 0.875
$ws.Cells.Item(117, 28).Value = 0.8500000000000001

# Row 118
$ws.Cells.Item(118, 2).Value = 5265465
$ws.Cells.Item(118, 6).Value = "Feirense"
$ws.Cells.Item(118, 7).Value = "Penafiel"
$ws.Cells.Item(118, 8).Value = 2
$ws.Cells.Item(118, 9).Value = 1
$ws.Cells.Item(118, 11).Value = 2.25
$ws.Cells.Item(118, 12).Value = 3.1
$ws.Cells.Item(118, 13).Value = 3.1
$ws.Cells.Item(118, 14).Value = 2.55
$ws.Cells.Item(118, 15).Value = 3.1
$ws.Cells.Item(118, 16).Value = 2.7
$ws.Cells.Item(118, 17).Value = 0
$ws.Cells.Item(118, 18).Value = 1.85
$ws.Cells.Item(118, 19).Value = 1.95
$ws.Cells.Item(118, 20).Value = 2.25
$ws.Cells.Item(118, 21).Value = 2.025
$ws.Cells.Item(118, 22).Value = 1.775
$ws.Cells.Item(118, 23).Value = 1.55
$ws.Cells.Item(118, 26).Value = 0.8500000000000001
$ws.Cells.Item(118, 28).Value = 1.025

# Row 185
$ws.Cells.Item(185, 2).Value = 6893191
$ws.Cells.Item(185, 6).Value = "Academico Viseu"
$ws.Cells.Item(185, 7).Value = "Vilaverdense"
$ws.Cells.Item(185, 11).Value = 1.833
$ws.Cells.Item(185, 12).Value = 3.4
$ws.Cells.Item(185, 14).Value = 1.65
$ws.Cells.Item(185, 15).Value = 3.75
$ws.Cells.Item(185, 16).Value = 4.2
$ws.Cells.Item(185, 17).Value = -0.75
$ws.Cells.Item(185, 18).Value = 1.95
$ws.Cells.Item(185, 19).Value = 1.85
$ws.Cells.Item(185, 20).Value = 2.5
$ws.Cells.Item(185, 21).Value = 1.975
$ws.Cells.Item(185, 22).Value = 1.825
$ws.Cells.Item(185, 24).Value = 2.75
$ws.Cells.Item(185, 26).Value = -1
$ws.Cells.Item(185, 27).Value = 0.8500000000000001
$ws.Cells.Item(185, 28).Value = -1
$ws.Cells.Item(185, 29).Value = 0.825

# Row 186
$ws.Cells.Item(186, 2).Value = 6893285
$ws.Cells.Item(186, 6).Value = "FC Porto B"
$ws.Cells.Item(186, 7).Value = "Tondela"
$ws.Cells.Item(186, 11).Value = 1.8
$ws.Cells.Item(186, 12).Value = 3.5
$ws.Cells.Item(186, 14).Value = 2.2
$ws.Cells.Item(186, 15).Value = 3.3
$ws.Cells.Item(186, 16).Value = 2.75
$ws.Cells.Item(186, 17).Value = -0.25
$ws.Cells.Item(186, 18).Value = 2
$ws.Cells.Item(186, 19).Value = 1.8
$ws.Cells.Item(186, 20).Value = 2.25
$ws.Cells.Item(186, 21).Value = 1.8
$ws.Cells.Item(186, 22).Value = 2
$ws.Cells.Item(186, 24).Value = 2.3
$ws.Cells.Item(186, 26).Value = -0.5
$ws.Cells.Item(186, 27).Value = 0.4
$ws.Cells.Item(186, 28).Value = -0.5
$ws.Cells.Item(186, 29).Value = 0.5

# Row 254
$ws.Cells.Item(254, 2).Value = 6893574
$ws.Cells.Item(254, 6).Value = "FC Porto B"
$ws.Cells.Item(254, 7).Value = "Feirense"
$ws.Cells.Item(254, 8).Value = 2
$ws.Cells.Item(254, 9).Value = 0
$ws.Cells.Item(254, 10).Value = "H"
$ws.Cells.Item(254, 11).Value = 2
$ws.Cells.Item(254, 12).Value = 3.3
$ws.Cells.Item(254, 13).Value = 3.2
$ws.Cells.Item(254, 14).Value = 1.727
$ws.Cells.Item(254, 15).Value = 3.6
$ws.Cells.Item(254, 18).Value = 1.775
$ws.Cells.Item(254, 19).Value = 2.025
$ws.Cells.Item(254, 21).Value = 1.825
$ws.Cells.Item(254, 22).Value = 1.975
$ws.Cells.Item(254, 23).Value = 0.7270000000000001
$ws.Cells.Item(254, 25).Value = -1
$ws.Cells.Item(254, 26).Value = 0.7749999999999999
$ws.Cells.Item(254, 27).Value = -1
$ws.Cells.Item(254, 29).Value = 0.9750000000000001

# Row 255
$ws.Cells.Item(255, 2).Value = 6893126
$ws.Cells.Item(255, 6).Value = "CD Mafra"
$ws.Cells.Item(255, 7).Value = "Leixoes"
$ws.Cells.Item(255, 8).Value = 0
$ws.Cells.Item(255, 9).Value = 1
$ws.Cells.Item(255, 10).Value = "A"
$ws.Cells.Item(255, 11).Value = 1.8
$ws.Cells.Item(255, 12).Value = 3.4
$ws.Cells.Item(255, 13).Value = 3.8
$ws.Cells.Item(255, 14).Value = 1.8
$ws.Cells.Item(255, 15).Value = 3.3
$ws.Cells.Item(255, 18).Value = 1.925
$ws.Cells.Item(255, 19).Value = 1.925
$ws.Cells.Item(255, 21).Value = 1.975
$ws.Cells.Item(255, 22).Value = 1.825
$ws.Cells.Item(255, 23).Value = -1
$ws.Cells.Item(255, 25).Value = 2.8
$ws.Cells.Item(255, 26).Value = -1
$ws.Cells.Item(255, 27).Value = 0.925
$ws.Cells.Item(255, 29).Value = 0.825

# Row 370
$ws.Cells.Item(370, 2).Value = 6893601
$ws.Cells.Item(370, 5).Value = 45339.45833333334
$ws.Cells.Item(370, 7).Value = "Benfica B"
$ws.Cells.Item(370, 11).Value = 2.45
$ws.Cells.Item(370, 13).Value = 2.875
$ws.Cells.Item(370, 14).Value = 2.3
$ws.Cells.Item(370, 16).Value = 3
$ws.Cells.Item(370, 18).Value = 2.1
$ws.Cells.Item(370, 19).Value = 1.775
$ws.Cells.Item(370, 21).Value = 1.925
$ws.Cells.Item(370, 22).Value = 1.925

# Row 371
$ws.Cells.Item(371, 2).Value = 6899171
$ws.Cells.Item(371, 5).Value = 45339.52083333334
$ws.Cells.Item(371, 6).Value = "Maritimo"
$ws.Cells.Item(371, 7).Value = "FC Porto B"
$ws.Cells.Item(371, 11).Value = 1.85
$ws.Cells.Item(371, 12).Value = 3.5
$ws.Cells.Item(371, 13).Value = 3.75
$ws.Cells.Item(371, 14).Value = 1.8
$ws.Cells.Item(371, 15).Value = 3.5
$ws.Cells.Item(371, 16).Value = 3.8
$ws.Cells.Item(371, 17).Value = -0.5
$ws.Cells.Item(371, 18).Value = 1.875
$ws.Cells.Item(371, 19).Value = 1.975
$ws.Cells.Item(371, 20).Value = 2.5
$ws.Cells.Item(371, 21).Value = 1.95
$ws.Cells.Item(371, 22).Value = 1.9

# Row 372
$ws.Cells.Item(372, 2).Value = 6899170
$ws.Cells.Item(372, 5).Value = 45339.625
$ws.Cells.Item(372, 6).Value = "UD Leiria"
$ws.Cells.Item(372, 7).Value = "CF Os Belenenses"
$ws.Cells.Item(372, 11).Value = 1.571
$ws.Cells.Item(372, 13).Value = 5.25
$ws.Cells.Item(372, 14).Value = 1.55
$ws.Cells.Item(372, 15).Value = 3.6
$ws.Cells.Item(372, 16).Value = 5.25
$ws.Cells.Item(372, 17).Value = -0.75
$ws.Cells.Item(372, 18).Value = 1.775
$ws.Cells.Item(372, 19).Value = 2.1
$ws.Cells.Item(372, 21).Value = 2.05
$ws.Cells.Item(372, 22).Value = 1.8

# Row 373
$ws.Cells.Item(373, 2).Value = 6893202
$ws.Cells.Item(373, 5).Value = 45340.33333333334
$ws.Cells.Item(373, 6).Value = "Academico Viseu"
$ws.Cells.Item(373, 7).Value = "UD Oliveirense"
$ws.Cells.Item(373, 11).Value = 1.727
$ws.Cells.Item(373, 13).Value = 4.5
$ws.Cells.Item(373, 14).Value = 1.65
$ws.Cells.Item(373, 16).Value = 4.75
$ws.Cells.Item(373, 18).Value = 1.875
$ws.Cells.Item(373, 19).Value = 1.975
$ws.Cells.Item(373, 20).Value = 2.25
$ws.Cells.Item(373, 21).Value = 1.825
$ws.Cells.Item(373, 22).Value = 2.025

# Row 374
$ws.Cells.Item(374, 2).Value = 6893164
$ws.Cells.Item(374, 5).Value = 45340.45833333334
$ws.Cells.Item(374, 6).Value = "Penafiel"
$ws.Cells.Item(374, 7).Value = "Santa Clara"
$ws.Cells.Item(374, 11).Value = 3.4
$ws.Cells.Item(374, 12).Value = 2.875
$ws.Cells.Item(374, 13).Value = 2.2
$ws.Cells.Item(374, 14).Value = 4
$ws.Cells.Item(374, 15).Value = 3
$ws.Cells.Item(374, 16).Value = 1.95
$ws.Cells.Item(374, 17).Value = 0.5
$ws.Cells.Item(374, 18).Value = 1.8
$ws.Cells.Item(374, 19).Value = 2.05
$ws.Cells.Item(374, 20).Value = 2
$ws.Cells.Item(374, 21).Value = 1.975
$ws.Cells.Item(374, 22).Value = 1.875

# Row 375
$ws.Cells.Item(375, 2).Value = 6893165
$ws.Cells.Item(375, 5).Value = 45340.52083333334
$ws.Cells.Item(375, 6).Value = "Feirense"
$ws.Cells.Item(375, 7).Value = "Nacional"
$ws.Cells.Item(375, 11).Value = 3.2
$ws.Cells.Item(375, 12).Value = 3.25
$ws.Cells.Item(375, 13).Value = 2.15
$ws.Cells.Item(375, 14).Value = 2.7
$ws.Cells.Item(375, 15).Value = 3.2
$ws.Cells.Item(375, 16).Value = 2.45
$ws.Cells.Item(375, 17).Value = 0
$ws.Cells.Item(375, 18).Value = 2
$ws.Cells.Item(375, 19).Value = 1.85
$ws.Cells.Item(375, 20).Value = 2.5
$ws.Cells.Item(375, 21).Value = 2.025
$ws.Cells.Item(375, 22).Value = 1.825

# Row 376
$ws.Cells.Item(376, 2).Value = 6899276
$ws.Cells.Item(376, 5).Value = 45341.625
$ws.Cells.Item(376, 6).Value = "Vilaverdense"
$ws.Cells.Item(376, 7).Value = "AVS"
$ws.Cells.Item(376, 11).Value = 4
$ws.Cells.Item(376, 12).Value = 3.4
$ws.Cells.Item(376, 13).Value = 1.8
$ws.Cells.Item(376, 14).Value = 3.8
$ws.Cells.Item(376, 15).Value = 3.4
$ws.Cells.Item(376, 16).Value = 1.833
$ws.Cells.Item(376, 17).Value = 0.5
$ws.Cells.Item(376, 18).Value = 1.95
$ws.Cells.Item(376, 19).Value = 1.9
$ws.Cells.Item(376, 20).Value = 2.25
$ws.Cells.Item(376, 21).Value = 1.8
$ws.Cells.Item(376, 22).Value = 2.05

# Remove the last row (row 377), which is being dropped entirely
$ws.Rows.Item(377).Delete()

$wb.Save()